$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5464854836463928
$ws.Range("B1").Value = 2.048335790634155
$ws.Range("D1").Value = 2.252710103988647
$ws.Range("E1").Value = 1.136104822158813
